# Automatic update (mar abr  6 17:33:39 CEST 2021)
# The "porcentaje-participacion" metric row was dropped from the source
# dataset, so the column that used to describe it now carries the
# "participacion" metric's identifiers instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "participacion"
$ws.Range("L3").Value = "iaest-measure:participacion"
